$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - F column updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3442
$ws1.Range("F4").Value = 135
$ws1.Range("F5").Value = 6982
$ws1.Range("F6").Value = 2493
$ws1.Range("F8").Value = 114
$ws1.Range("F12").Value = 35
$ws1.Range("F13").Value = 175
$ws1.Range("F14").Value = 578
$ws1.Range("F15").Value = 44

# Sheet "全部类型" (All types) - F column updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3442
$ws4.Range("F5").Value = 135
$ws4.Range("F6").Value = 6982
$ws4.Range("F7").Value = 2493
$ws4.Range("F9").Value = 114
$ws4.Range("F13").Value = 35
$ws4.Range("F14").Value = 175
$ws4.Range("F15").Value = 578
$ws4.Range("F16").Value = 44
